$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N2").Value = 1.1
$ws.Range("R2").Value = 1.24
$ws.Range("S2").Value = 2.2
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("G3").Value = 2.66
$ws.Range("L3").Value = 1.33
$ws.Range("N3").Value = 2.6
$ws.Range("S3").Value = 2.56
$ws.Range("W3").Value = 1.6
$ws.Range("G4").Value = 2.94
$ws.Range("S4").Value = 4.2
$ws.Range("K6").Value = 3.6
$ws.Range("N7").Value = 1.98
$ws.Range("P7").Value = 1.98
$ws.Range("N8").Value = 1.89
$ws.Range("G9").Value = 3.3
$ws.Range("W9").Value = 1.38
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 5.3
$ws.Range("K10").Value = 4
$ws.Range("I11").Value = 3.05
$ws.Range("K11").Value = 3.65
$ws.Range("V11").Value = 1.49
$ws.Range("F12").Value = 3.8
$ws.Range("R12").Value = 1.25
$ws.Range("U12").Value = 1.91
$ws.Range("G13").Value = 2.46
$ws.Range("H13").Value = 3.1
$ws.Range("W13").Value = 1.69
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 3.9
$ws.Range("J14").Value = 3.25
$ws.Range("G16").Value = 2.04
$ws.Range("W16").Value = 1.96
$ws.Range("J17").Value = 4.3
$ws.Range("I18").Value = 4.1
$ws.Range("J18").Value = 3.8
$ws.Range("N18").Value = 4.1
$ws.Range("R18").Value = 1.51
$ws.Range("S18").Value = 2.6
$ws.Range("T18").Value = 1.6
$ws.Range("U18").Value = 2.38
$ws.Range("V18").Value = 1.33
$ws.Range("X18").Value = 26
$ws.Range("Y18").Value = 22
$ws.Range("Z18").Value = 36
$ws.Range("AA18").Value = 80
$ws.Range("AB18").Value = 15
$ws.Range("AC18").Value = 11.5
$ws.Range("AD18").Value = 19
$ws.Range("AE18").Value = 48
$ws.Range("AF18").Value = 18.5
$ws.Range("AG18").Value = 13.5
$ws.Range("AH18").Value = 19.5
$ws.Range("AI18").Value = 50
$ws.Range("AJ18").Value = 30
$ws.Range("AK18").Value = 24
$ws.Range("AL18").Value = 36
$ws.Range("AM18").Value = 80
$ws.Range("AN18").Value = 13.5
$ws.Range("AO18").Value = 36
$ws.Range("F20").Value = 2.72
$ws.Range("K20").Value = 3.2
$ws.Range("N21").Value = 1.53
$ws.Range("P21").Value = 1.53
$ws.Range("H22").Value = 2.08
$ws.Range("I22").Value = 2.18
$ws.Range("R22").Value = 1.22
$ws.Range("U22").Value = 1.77
$ws.Range("V22").Value = 1.84
$ws.Range("AM23").Value = 150
$ws.Range("G24").Value = 2.52
$ws.Range("N24").Value = 1.25
$ws.Range("P24").Value = 1.25
$ws.Range("R24").Value = 1.12
$ws.Range("M25").Value = 1.06
$ws.Range("T25").Value = 1.77
$ws.Range("W25").Value = 2.2
$ws.Range("Y25").Value = 19.5
$ws.Range("AH25").Value = 18.5
$ws.Range("AN25").Value = 9.6
$ws.Range("P26").Value = 1.78
$ws.Range("T26").Value = 1.97
$ws.Range("V26").Value = 1.3
$ws.Range("W26").Value = 1.89
$ws.Range("AM26").Value = 130
$ws.Range("G27").Value = 1.74
$ws.Range("L27").Value = 1.49
$ws.Range("M27").Value = 1.08
$ws.Range("N27").Value = 1.73
$ws.Range("O27").Value = 1.4
$ws.Range("R27").Value = 1.22
$ws.Range("S27").Value = 3.55
$ws.Range("T27").Value = 1.92
$ws.Range("U27").Value = 1.64
$ws.Range("V27").Value = 1.16
$ws.Range("W27").Value = 2.34
$ws.Range("X27").Value = 1000
$ws.Range("Y27").Value = 25
$ws.Range("Z27").Value = 75
$ws.Range("AA27").Value = 1000
$ws.Range("AB27").Value = 9.6
$ws.Range("AC27").Value = 12.5
$ws.Range("AD27").Value = 30
$ws.Range("AE27").Value = 1000
$ws.Range("AF27").Value = 12.5
$ws.Range("AG27").Value = 12
$ws.Range("AH27").Value = 36
$ws.Range("AI27").Value = 1000
$ws.Range("AJ27").Value = 24
$ws.Range("AK27").Value = 24
$ws.Range("AL27").Value = 70
$ws.Range("AM27").Value = 1000
$ws.Range("AN27").Value = 1000
$ws.Range("AO27").Value = 1000
$ws.Range("L28").Value = 1.5
$ws.Range("M28").Value = 1.01
$ws.Range("N28").Value = 1.63
$ws.Range("O28").Value = 1.44
$ws.Range("Q28").Value = 2.12
$ws.Range("R28").Value = 1.18
$ws.Range("S28").Value = 3.85
$ws.Range("T28").Value = 1.64
$ws.Range("U28").Value = 1.61
$ws.Range("V28").Value = 1.47
$ws.Range("W28").Value = 1.45
$ws.Range("X28").Value = 14.5
$ws.Range("Y28").Value = 13.5
$ws.Range("Z28").Value = 26
$ws.Range("AA28").Value = 70
$ws.Range("AB28").Value = 13.5
$ws.Range("AC28").Value = 10
$ws.Range("AD28").Value = 19
$ws.Range("AE28").Value = 55
$ws.Range("AF28").Value = 28
$ws.Range("AG28").Value = 19
$ws.Range("AH28").Value = 980
$ws.Range("AI28").Value = 85
$ws.Range("AJ28").Value = 75
$ws.Range("AK28").Value = 55
$ws.Range("AL28").Value = 85
$ws.Range("AM28").Value = 1000
$ws.Range("AN28").Value = 1000
$ws.Range("AO28").Value = 1000
$ws.Range("G29").Value = 2.36
$ws.Range("I29").Value = 4.5
$ws.Range("W29").Value = 1.73
